# "need to fix up task 2"
# Applies the documented edits to the Custom Physics Simulation write-up:
#   1-3, 5) a few phrases get a grammar-style "is/are able to" split out
#      (net visible wording is unchanged, only how it's run-split)
#   4) the old "second improvement" (efficiency) paragraph under
#      "3.2 - Improvement #2" is removed, along with the blank paragraph
#      that followed it, and replaced with a new closing paragraph about
#      angular velocity / torque that sits right before the page break.

$d = $word.ActiveDocument
$rsquo = [char]0x2019

function Touch-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Edit 1: "...that I've created is able to successfully demonstrate..." --
Touch-Text "is able to successfully " "is able to successfully "

# --- Edit 2: "...physics engine I've created is able to showcase..." -------
Touch-Text ("I" + $rsquo + "ve created is able to showcase the effect that collision has on ") `
           ("I" + $rsquo + "ve created is able to showcase the effect that collision has on ")

# --- Edit 3: "The physic bodies are able to interact together..." ----------
Touch-Text "are able to interact together through the use of collision between each object. Objects that are dynamic are able to move due to factors such as" `
           "are able to interact together through the use of collision between each object. Objects that are dynamic are able to move due to factors such as"

# --- Edit 5: "The third party libraries that were used..." -----------------
Touch-Text "The third party libraries that were used within this custom physics engine were things such as OpenGL and GLM (OpenGL Mathematics). " `
           "The third party libraries that were used within this custom physics engine were things such as OpenGL and GLM (OpenGL Mathematics). "

# --- Edit 4: drop the old Improvement #2 body + blank line, add the new ----
#             "Another improvement ... angular velocity ..." paragraph text
#             right before the section's page break.
$target = $null
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*second improvement that could be made to the Physics Engine*") {
        $target = $i
    }
}
if ($target -ne $null) {
    # Remove the "the second improvement ... streamlined." paragraph entirely.
    $d.Paragraphs($target).Range.Delete()
    # Remove the now-following blank paragraph too.
    $d.Paragraphs($target).Range.Delete()
    # $target now refers to the paragraph holding the page break; insert the
    # new text as its own run ahead of the <w:br w:type="page"/>.
    $pageBreakPara = $d.Paragraphs($target)
    $point = $d.Range($pageBreakPara.Range.Start, $pageBreakPara.Range.Start)
    $point.InsertAfter("Another improvement that could be made to the physics engine is to incorporate better usage of angular velocity and torque. Angular momentum refers to the rotational counterpart of linear momentum. The addition to angular velocity in the physics engine would allow for the circles to react more normal like. ")
}

Write-Host "done"
